$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24-30: mark as completed ("Y") instead of "N"
$ws.Range("D24:D30").Value = "Y"

# The trailing block of rows (149-170) is no longer part of this roster;
# clear it out (they sit at the very bottom of the used range, so the
# sheet's dimension/used-range shrinks back down to row 148).
$ws.Range("B149:D170").ClearContents()

# Column C (raw minute counts) is no longer shown on screen.
$ws.Range("C1").EntireColumn.Hidden = $true

# Update view/selection to where the user left off.
[void]$ws.Range("H17").Select()
